$wb = $excel.ActiveWorkbook

# --- Add the new "attacks" sheet at the end of the tab order ---
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "attacks"
[void]$ws.Cells.Select()

# --- Populate the attack/ability data, column by column (matches authoring order) ---
$data = @(
  @("id",             "name",          "desc",                                              "tags"),
  @("attack-claw1",   "Claw Attack",   "swipe{!s} at {t} with their sharp claws",            "ability,action,attack,melee,physical"),
  @("attack-tongue1", "Tongue Smack",  "smacks{!s} {t} with their tongue",                   "ability,action,attack,ranged,physical"),
  @("attack-spit1",   "Acid Spit",     "spit{!s} acid at {t}",                               "ability,action,attack,ranged,acid")
)
$numData = @(
  @("dmg-max", "dmg-min", "spd"),
  @(4, 1, 80),
  @(4, 1, 80),
  @(8, 4, 120)
)

for ($c = 0; $c -lt 4; $c++) {
  for ($r = 0; $r -lt 4; $r++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
  }
}
for ($c = 0; $c -lt 3; $c++) {
  for ($r = 0; $r -lt 4; $r++) {
    $ws.Cells.Item($r + 1, $c + 5).Value = $numData[$r][$c]
  }
}

# Header row formatting: reuse the same wrap/vertical-center style used by the
# other sheets' header rows, and match their 30pt row height.
$headerSrc = $wb.Worksheets.Item("weapons").Range("A1:G1")
$headerSrc.Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)
$ws.Rows.Item(1).RowHeight = 30

# Approximate the bespoke column widths used on the new sheet.
$ws.Columns.Item(1).ColumnWidth = 13.6
$ws.Columns.Item(2).ColumnWidth = 12.9
$ws.Columns.Item(3).ColumnWidth = 33.7
$ws.Columns.Item(4).ColumnWidth = 32.7
$ws.Columns.Item(5).ColumnWidth = 4.7
$ws.Columns.Item(6).ColumnWidth = 4.7
$ws.Columns.Item(7).ColumnWidth = 3.3

# --- Restore/adjust the selections left on the pre-existing sheets ---
$wsMonsters = $wb.Worksheets.Item("monsters")
[void]$wsMonsters.Activate()
[void]$wsMonsters.Range("I8").Select()

$wsWeapons = $wb.Worksheets.Item("weapons")
[void]$wsWeapons.Activate()
[void]$wsWeapons.Range("G1:J1").Select()

$wsArmor = $wb.Worksheets.Item("armor")
[void]$wsArmor.Activate()
[void]$wsArmor.Range("G6").Select()

# --- The newly added "attacks" sheet is the one left active ---
[void]$ws.Activate()

Write-Output "done"
